$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '43.450.75'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.02%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.232.76'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.10%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '258.15'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.28%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '79.13'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +5.12%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.620'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.04%  '
$ws.Range("E8").Value = '  -0.10%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.600'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.32%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '43.22'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +5.17%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0923'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("E13").Value = '  +0.85%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '2.567.38'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.03%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '14.60'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.44%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.234.09'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.17%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.794'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.49%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '43.349.17'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +0.82%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '71.34'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("E22").Value = '  +5.71%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '230.10'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -0.14%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '41.97'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +6.80%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.84'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +2.46%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '172.88'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.58%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '20.50'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.48%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0869'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +9.22%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  +0.94%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0370'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +13.40%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.48'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("E38").Value = '  -4.32%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '13.22'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +7.56%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.87'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +18.10%  '
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("E42").Value = '  -0.34%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '61.65'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.38%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.36'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '103.61'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '8.57'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.83%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.471'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.40%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0982'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("E50").Value = '  +1.29%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.47'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +23.17%  '
